$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values are written as text, even though many look
# like numbers (e.g. "1.00", "10.30") - force text number format first so
# Excel does not coerce them into numeric values and lose formatting.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.492.83"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.554.57"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.16"
$ws.Range("E5").Value = "  +1.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.55"
$ws.Range("E6").Value = "  +3.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.553.33"
$ws.Range("E7").Value = "  +3.11%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.492"
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("E10").Value = "  +2.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.05"
$ws.Range("E11").Value = "  -4.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.395"
$ws.Range("E12").Value = "  +4.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.160.37"
$ws.Range("E13").Value = "  +3.20%  "
$ws.Range("E14").Value = "  +2.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.24"
$ws.Range("E15").Value = "  +1.98%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.556.36"
$ws.Range("E16").Value = "  +3.24%  "
$ws.Range("E17").Value = "  +1.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.535.11"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.30"
$ws.Range("E19").Value = "  +3.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.92"
$ws.Range("E20").Value = "  +1.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.30"
$ws.Range("E21").Value = "  +4.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "395.68"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.573"
$ws.Range("E23").Value = "  +4.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.700.63"
$ws.Range("E24").Value = "  +3.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.19"
$ws.Range("E25").Value = "  +1.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  +9.70%  "
$ws.Range("E28").Value = "  +9.75%  "
$ws.Range("E29").Value = "  +0.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.30"
$ws.Range("E30").Value = "  +1.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.44"
$ws.Range("E31").Value = "  +2.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.568.57"
$ws.Range("E32").Value = "  +3.29%  "
$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.148"
$ws.Range("E34").Value = "  +0.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.80"
$ws.Range("E35").Value = "  +3.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.30"
$ws.Range("E36").Value = "  +6.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.08"
$ws.Range("E37").Value = "  +2.35%  "
$ws.Range("E38").Value = "  +3.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "167.87"
$ws.Range("E39").Value = "  -2.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.05"
$ws.Range("E40").Value = "  +4.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0814"
$ws.Range("E41").Value = "  +4.94%  "
$ws.Range("E42").Value = "  +1.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.66"
$ws.Range("E43").Value = "  +15.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "43.06"
$ws.Range("E44").Value = "  -0.98%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("E46").Value = "  +0.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.71"
$ws.Range("E47").Value = "  +5.06%  "
$ws.Range("E48").Value = "  +8.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.480.46"
$ws.Range("E49").Value = "  +12.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.83"
$ws.Range("E50").Value = "  +3.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.35"
$ws.Range("E51").Value = "  +18.92%  "
